$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Add the "C" validation column (order chosen to reproduce the shared-string table order)
$ws.Range("C2").Value = "Invalid"
$ws.Range("C3").Value = "Invalid"
$ws.Range("C4").Value = "Invalid"
$ws.Range("C5").Value = "Valid"
$ws.Range("C1").Value = "Expected "

# 2. Fix email typo in B2 (text only; underlying hyperlink target is left unchanged)
$ws.Range("B2").Value = "sw@g.com"

# 3. Build a bold + yellow-fill header style on a scratch cell, then copy/paste the
#    formatting (as a single operation) onto the header row so only one new cell
#    style (cellXf) is produced instead of one per property.
$scratch = $ws.Range("Z1")
$scratch.Font.Bold = $true
$scratch.Interior.Color = 65535

$header = $ws.Range("A1:C1")
$scratch.Copy()
$header.PasteSpecial(-4122)  # xlPasteFormats
$scratch.Clear()

# 4. Update the selected cell shown when the workbook is opened
$ws.Range("B2").Select()

# 5. Page setup (portrait orientation)
$ws.PageSetup.Orientation = 1
